$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend row-number column (A) down to rows 15-18, copying the
#     existing numeric style (border + bold + centered) from A14 so the
#     new cells pick up the same cellXf (s="1") as the rest of column A.
$ws.Range("A14").Copy($ws.Range("A15:A18")) | Out-Null
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16

# --- Column B (Buying Opportunity) ---
$ws.Range("B2").Value = "NSE:AXISGOLD"
$ws.Range("B3").Value = "NSE:BSLGOLDETF"
$ws.Range("B4").Value = "NSE:EGOLD"
$ws.Range("B5").Value = "NSE:GOLDBEES"
$ws.Range("B6").Value = "NSE:GOLDSHARE"
$ws.Range("B7").Value = "NSE:LTGILTBEES"
$ws.Range("B8").Value = "NSE:NMDC"

# --- Column C (support Zone) ---
$ws.Range("C2").Value = "NSE:AJANTPHARM"
$ws.Range("C3").Value = "NSE:APOLLOTYRE"
$ws.Range("C4").Value = "NSE:AUBANK"
$ws.Range("C5").Value = "NSE:BDL"
$ws.Range("C6").Value = "NSE:BLUESTARCO"
$ws.Range("C7").Value = "NSE:COCHINSHIP"
$ws.Range("C8").Value = "NSE:CROMPTON"
$ws.Range("C9").Value = "NSE:GODREJAGRO"
$ws.Range("C10").Value = "NSE:GRANULES"
$ws.Range("C11").Value = "NSE:HEG"
$ws.Range("C12").Value = "NSE:JSWENERGY"
$ws.Range("C13").Value = "NSE:KAJARIACER"
$ws.Range("C14").Value = "NSE:MSUMI"
$ws.Range("C15").Value = "NSE:NAZARA"
$ws.Range("C16").Value = "NSE:OFSS"
$ws.Range("C17").Value = "NSE:REDTAPE"
$ws.Range("C18").Value = "NSE:RVNL"

# --- Column E (Short buildup) ---
$ws.Range("E2").Value = "NSE:BIOCON"
$ws.Range("E3").Value = "NSE:DABUR"
$ws.Range("E4").Value = "NSE:ESCORTS"

# --- Column F2 (FII ENTERING) is cleared - no longer populated ---
$ws.Range("F2").ClearContents() | Out-Null
